# Atualização de bases das ligas, do dia: 04-04-2024 às 23:22
#
# The source data rows got re-sorted (by kickoff id) and this re-ordering
# shows up in the OOXML diff as full-row content swaps: column A (the
# running row index) stays put, but everything from column B (id) through
# column AC (PL_AhUnder) moves between rows.
#
# Below: simple pairwise swaps, plus one 3-way rotation (rows 208/209/210).
# NOTE: use positional parameters for the helper functions - named
# parameters (-RowA 194) do not bind correctly in this PowerShell host.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($RowA, $RowB) {
    $rangeA = $ws.Range("B" + $RowA + ":AC" + $RowA)
    $rangeB = $ws.Range("B" + $RowB + ":AC" + $RowB)
    $valsA = $rangeA.Value2
    $valsB = $rangeB.Value2
    $rangeA.Value = $valsB
    $rangeB.Value = $valsA
}

function Rotate3-RowData($RowA, $RowB, $RowC) {
    # RowA <- old RowC, RowB <- old RowA, RowC <- old RowB
    $rangeA = $ws.Range("B" + $RowA + ":AC" + $RowA)
    $rangeB = $ws.Range("B" + $RowB + ":AC" + $RowB)
    $rangeC = $ws.Range("B" + $RowC + ":AC" + $RowC)
    $valsA = $rangeA.Value2
    $valsB = $rangeB.Value2
    $valsC = $rangeC.Value2
    $rangeA.Value = $valsC
    $rangeB.Value = $valsA
    $rangeC.Value = $valsB
}

Swap-RowData 194 195
Rotate3-RowData 208 209 210
Swap-RowData 214 215
Swap-RowData 221 222
Swap-RowData 227 228
Swap-RowData 248 249
Swap-RowData 258 259
Swap-RowData 264 265
